# Summer 24 week 3 inputs
# Update a handful of matchup-average cells on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.41
$ws.Range("E3").Value = 1.34
$ws.Range("C5").Value = 1.32
$ws.Range("F6").Value = 1.19
